$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$testData = "Username:standard_user" + [char]10 + "Password:secret_sauce"

for ($r = 9; $r -le 30; $r++) {
    $ws.Range("G$r").Value = $testData
}

# Update the sheet view state (zoom/scroll/selection) to match the saved file
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 55
$ws.Range("G19").Select()
